# Update four flowchart boxes whose text is currently "item" to "Item changed".
# These four shapes (identified by their stable p:cNvPr @id values) are the
# first three "item" boxes in document order (ids 4, 7, 8) plus the box with
# id 5 ("Flowchart: Process 4") further down the shape tree.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$targetIds = @(4, 7, 8, 5)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($targetIds -contains $shp.Id) {
        if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "item") {
                $shp.TextFrame.TextRange.Text = "Item changed"
            }
        }
    }
}
